# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Frambuesa" (Raspberry) at
# Mercado Mayorista Lo Valledor de Santiago, pushing the existing rows
# 167-188 down to 168-189.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 167 (shifts rows 167-188 down to 168-189,
# and copies formatting - e.g. the date number format on column D - from
# the row being pushed down).
$ws.Rows.Item(167).Insert()

# Populate the new row 167 with the new weekly record.
$ws.Cells.Item(167, 1).Value = 6
$ws.Cells.Item(167, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(167, 3).Value = "Metropolitana"
$ws.Cells.Item(167, 4).Value = 44644
$ws.Cells.Item(167, 5).Value = 13
$ws.Cells.Item(167, 6).Value = "Fruta"
$ws.Cells.Item(167, 7).Value = 100101
$ws.Cells.Item(167, 8).Value = "Berries"
$ws.Cells.Item(167, 9).Value = 100101004
$ws.Cells.Item(167, 10).Value = "Frambuesa"
$ws.Cells.Item(167, 11).Value = "Sin especificar"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 170
$ws.Cells.Item(167, 14).Value = 8000
$ws.Cells.Item(167, 15).Value = 8000
$ws.Cells.Item(167, 16).Value = 8000
$ws.Cells.Item(167, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(167, 18).Value = "Provincia de Linares"
$ws.Cells.Item(167, 19).Value = 4000
$ws.Cells.Item(167, 20).Value = 2
